$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31 is "urban_nr" - incorporate new urban net returns data
$ws.Range("E31").Value = $null
$ws.Range("F31").Value = $null
$ws.Range("G31").Value = 9314.0947265625
$ws.Range("H31").Value = 8957.0625
$ws.Range("I31").Value = 8047.3134765625
